$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.05320733333333333
$ws.Range("H2").Value = 0.159622
$ws.Range("I2").Value = 0.05437835303993056
$ws.Range("J2").Value = 0.05437835303993056
$ws.Range("M2").Value = 0.09551033333333332
$ws.Range("N2").Value = 0.286531
$ws.Range("O2").Value = 0.0198020999427218
$ws.Range("P2").Value = 0.0198020999427218
$ws.Range("Q2").Value = 0.005081850142444443
$ws.Range("R2").Value = 0.045736651282
$ws.Range("S2").Value = 0.001076805581617315
$ws.Range("T2").Value = 0.001076805581617315
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.05320733333333333
$ws.Range("H3").Value = 0.159622
$ws.Range("I3").Value = 0.05437835303993056
$ws.Range("J3").Value = 0.05437835303993056
$ws.Range("O3").Value = 0.07175622098770619
$ws.Range("P3").Value = 0.07175622098770619
$ws.Range("Q3").Value = 0.01841493391622222
$ws.Range("R3").Value = 0.165734405246
$ws.Range("S3").Value = 0.003901985117680762
$ws.Range("T3").Value = 0.003901985117680762
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.05320733333333333
$ws.Range("H4").Value = 0.159622
$ws.Range("I4").Value = 0.05437835303993056
$ws.Range("J4").Value = 0.05437835303993056
$ws.Range("M4").Value = 4.381634666666667
$ws.Range("N4").Value = 13.144904
$ws.Range("O4").Value = 0.9084416790695721
$ws.Range("P4").Value = 0.9084416790695721
$ws.Range("Q4").Value = 0.2331350962542222
$ws.Range("R4").Value = 2.098215866288
$ws.Range("S4").Value = 0.04939956234063248
$ws.Range("T4").Value = 0.04939956234063248
$ws.Range("I5").Value = 0.779193676083227
$ws.Range("J5").Value = 0.779193676083227
$ws.Range("M5").Value = 0.09551033333333332
$ws.Range("N5").Value = 0.286531
$ws.Range("O5").Value = 0.0198020999427218
$ws.Range("P5").Value = 0.0198020999427218
$ws.Range("Q5").Value = 0.072818415278
$ws.Range("R5").Value = 0.655365737502
$ws.Range("S5").Value = 0.01542967104853686
$ws.Range("T5").Value = 0.01542967104853686
$ws.Range("I6").Value = 0.779193676083227
$ws.Range("J6").Value = 0.779193676083227
$ws.Range("O6").Value = 0.07175622098770619
$ws.Range("P6").Value = 0.07175622098770619
$ws.Range("S6").Value = 0.05591199361325119
$ws.Range("T6").Value = 0.05591199361325119
$ws.Range("I7").Value = 0.779193676083227
$ws.Range("J7").Value = 0.779193676083227
$ws.Range("M7").Value = 4.381634666666667
$ws.Range("N7").Value = 13.144904
$ws.Range("O7").Value = 0.9084416790695721
$ws.Range("P7").Value = 0.9084416790695721
$ws.Range("Q7").Value = 3.340619612752
$ws.Range("R7").Value = 30.065576514768
$ws.Range("S7").Value = 0.7078520114214389
$ws.Range("T7").Value = 0.7078520114214389
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.162844
$ws.Range("H8").Value = 0.488532
$ws.Range("I8").Value = 0.1664279708768425
$ws.Range("J8").Value = 0.1664279708768425
$ws.Range("M8").Value = 0.09551033333333332
$ws.Range("N8").Value = 0.286531
$ws.Range("O8").Value = 0.0198020999427218
$ws.Range("P8").Value = 0.0198020999427218
$ws.Range("Q8").Value = 0.01555328472133333
$ws.Range("R8").Value = 0.139979562492
$ws.Range("S8").Value = 0.003295623312567628
$ws.Range("T8").Value = 0.003295623312567628
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.162844
$ws.Range("H9").Value = 0.488532
$ws.Range("I9").Value = 0.1664279708768425
$ws.Range("J9").Value = 0.1664279708768425
$ws.Range("O9").Value = 0.07175622098770619
$ws.Range("P9").Value = 0.07175622098770619
$ws.Range("Q9").Value = 0.05635992843066666
$ws.Range("R9").Value = 0.5072393558759999
$ws.Range("S9").Value = 0.01194224225677424
$ws.Range("T9").Value = 0.01194224225677424
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.162844
$ws.Range("H10").Value = 0.488532
$ws.Range("I10").Value = 0.1664279708768425
$ws.Range("J10").Value = 0.1664279708768425
$ws.Range("M10").Value = 4.381634666666667
$ws.Range("N10").Value = 13.144904
$ws.Range("O10").Value = 0.9084416790695721
$ws.Range("P10").Value = 0.9084416790695721
$ws.Range("Q10").Value = 0.7135229156586667
$ws.Range("R10").Value = 6.421706240928
$ws.Range("S10").Value = 0.1511901053075006
$ws.Range("T10").Value = 0.1511901053075006
